# Re-applies the "resolve and classify+summarise" steps after changes
# to the mapping file, for the Haryana SoIB_summaries workbook.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "Range Status": the per-category species counts all collapsed to
# zero (no species resolved to a range bucket this run), and the
# percentage column is no longer produced at all.
# ----------------------------------------------------------------------
$wsRange = $wb.Worksheets.Item("Range Status")

$wsRange.Range("B2:B7").Value = 0
$wsRange.Range("C2:C7").ClearContents()

# ----------------------------------------------------------------------
# Sheet "Species qualification": no species selected for the Range
# Analysis this run.
# ----------------------------------------------------------------------
$wsQual = $wb.Worksheets.Item("Species qualification")
$wsQual.Range("B5").Value = 0

# ----------------------------------------------------------------------
# Sheet "High Priority break-up": the "Range" category no longer
# contributes any high-priority species, so its row is dropped (shifting
# "IUCN" up into row 3), and the surviving rows get refreshed
# counts/percentages.
# ----------------------------------------------------------------------
$wsBreakup = $wb.Worksheets.Item("High Priority break-up")

# Drop the old "Range" row (row 3); "IUCN" (old row 4) shifts up to row 3.
$wsBreakup.Rows.Item(3).Delete()

$wsBreakup.Range("E2").Value = 7.7

$wsBreakup.Range("A3").Value = "IUCN"
$wsBreakup.Range("B3").Value = 12
$wsBreakup.Range("C3").Value = 92.3
$wsBreakup.Range("D3").Value = 12
$wsBreakup.Range("E3").Value = 92.3
